$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Bold the title paragraph ("Forespørsel om å delta i intervju...")
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.Bold = 1

# ------------------------------------------------------------------
# 2) Remove the old "_GoBack" bookmark that currently sits between
#    "Med vennlig" + " " and "hilsen", merging " " and "hilsen" into a
#    single " hilsen" run while keeping "Med vennlig" as its own run.
#    A temporary bookmark is used as a merge barrier so the text edit
#    does not fuse into the preceding "Med vennlig" run.
# ------------------------------------------------------------------
$mv = $d.Content
$mv.Find.Execute("Med vennlig", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$barrierPos = $mv.End
$barrierRange = $d.Range($barrierPos, $barrierPos)
$d.Bookmarks.Add("zzBarrier", $barrierRange) | Out-Null

$hilsenRange = $d.Content
$hilsenRange.Find.Execute("hilsen", $true, $false, $false, $false, $false, $true, 1, $false, "zzHILSENzz", 2) | Out-Null

$fixRange = $d.Content
$fixRange.Find.Execute("zzHILSENzz", $true, $false, $false, $false, $false, $true, 1, $false, "hilsen", 2) | Out-Null

$d.Bookmarks.Item("zzBarrier").Delete()
$d.Bookmarks.Item("_GoBack").Delete()

# ------------------------------------------------------------------
# 3) Split "...vil omfatte planlegging i tillegg..." so a fresh
#    "_GoBack" bookmark surrounds the "i" in "planlegging" (mirroring
#    Word's auto-tracked last-edit-location bookmark).
# ------------------------------------------------------------------
$planRange = $d.Content
$planRange.Find.Execute("planlegging", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$iStart = $planRange.Start + 8
$iRange = $d.Range($iStart, $iStart + 1)
$d.Bookmarks.Add("_GoBack", $iRange) | Out-Null
